$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row (row 8) that mirrors row 7 (same tag/values/formatting)
# but replaces the formula in column D with a hard-coded IF(5=6, ...) check.
$ws.Range("A7:F7").Copy()
$ws.Range("A8:F8").PasteSpecial(-4104)
$ws.Cells.Item(8, 4).Formula = '=IF(5=6,"presion no estable","presion estable")'

# Move/record the active selection to D9
$ws.Range("D9").Select()
